$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 35-39: timesheet entries -------------------------------------------------

# Row 35: Giovanni, 30/01/2018, 2h
$ws.Range("A35").Value = [DateTime]"2018-01-30"
$ws.Range("B35").Value = "Giovanni"
$ws.Range("C35").Value = "implementazione pulsante ricicla, potenzia da finire"
$ws.Range("D35").Value = 2 / 24
$ws.Rows.Item(35).RowHeight = 28.8

# Row 36: Mirko, 31/01/2018, 3h
$ws.Range("A36").Value = [DateTime]"2018-01-31"
$ws.Range("B36").Value = "Mirko"
$ws.Range("C36").Value = "implementazione metodi gerarchia"
$ws.Range("D36").Value = 3 / 24
$ws.Rows.Item(36).RowHeight = 28.8

# Row 37: Giovanni, 31/01/2018, 2h
$ws.Range("A37").Value = [DateTime]"2018-01-31"
$ws.Range("B37").Value = "Giovanni"
$ws.Range("C37").Value = "fine implementazione potenzia"
$ws.Range("D37").Value = 2 / 24
$ws.Rows.Item(37).RowHeight = 28.8

# Row 38: Mirko, 31/01/2018, 2h
$ws.Range("A38").Value = [DateTime]"2018-01-31"
$ws.Range("B38").Value = "Mirko"
$ws.Range("C38").Value = "crea e trasforma"
$ws.Range("D38").Value = 2 / 24

# Row 39: Giovanni, 01/02/2018, 1h
$ws.Range("A39").Value = [DateTime]"2018-02-01"
$ws.Range("B39").Value = "Giovanni"
$ws.Range("C39").Value = "implementazione crea"
$ws.Range("D39").Value = 1 / 24

# Row 40: only A40, date/time format with underline font, left blank.
$ws.Range("A40").NumberFormat = "dd/mm/yyyy\ hh:mm:ss"
$ws.Range("A40").Font.Underline = $true

# Update view: selection moves to the new bottom cell.
$ws.Range("A40").Select() | Out-Null
